$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) rows 2-11 from 45212 (2023-10-13) to 45221 (2023-10-22)
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45221
}
